# The deck's theme ("Integral") is swapped for the stock PowerPoint
# "Office Theme" palette. The two theme parts in the package
# (ppt/theme/theme1.xml - the live slide-master theme driving every
# slide's color scheme - and ppt/theme/theme2.xml, an orphaned part
# wired only to the notes master and not reachable through the
# PowerPoint object model) previously held "Integral" and "Office
# Theme" respectively; after the edit theme1.xml carries the Office
# colors. We reproduce that by rewriting the 12 theme color slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) that PowerPoint
# exposes via ThemeColorScheme, in index order, to the standard
# Office Theme RGB values.

function Get-BGR($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = Get-BGR 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = Get-BGR 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = Get-BGR 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = Get-BGR 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = Get-BGR 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = Get-BGR 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = Get-BGR 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = Get-BGR 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = Get-BGR 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = Get-BGR 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = Get-BGR 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = Get-BGR 0x95 0x4F 0x72   # folHlink 954F72
